$d = $word.ActiveDocument

# 1. Collapse the split "Part II of the / Advance Notice / of Methodological..."
#    runs (which were split apart by proofErr gramStart/gramEnd markers) into a
#    single run by doing a find/replace across the whole phrase. Word's
#    find/replace naturally merges the matched runs into one and drops the
#    now-irrelevant proofErr markers.
$titlePara = $d.Paragraphs(1).Range
$titlePara.Find.Execute(
    "Part II of the Advance Notice of Methodological Changes for CY 2019 for Medicare Advantage, Part D and 2019 draft Call Letter   ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Part II of the Advance Notice of Methodological Changes for CY 2019 for Medicare Advantage, Part D and 2019 draft Call Letter   ",
    2)

# 2. Move the "_GoBack" bookmark: it used to wrap nothing (zero-width, sitting
#    right after "March 2, 2018"); now it should span from the very start of
#    the document (start of the title paragraph) through to the very end of
#    the last paragraph's text.
$d.Bookmarks("_GoBack").Delete()
$start = $d.Paragraphs(1).Range.Start
$end = $d.Paragraphs($d.Paragraphs.Count).Range.End
$span = $d.Range($start, $end)
$d.Bookmarks.Add("_GoBack", $span)
